$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.046867666666667
$ws.Range("H2").Value = 3.140603
$ws.Range("I2").Value = 0.000687505225377314
$ws.Range("J2").Value = 0.000687505225377314
$ws.Range("M2").Value = 10.34761366666667
$ws.Range("N2").Value = 31.042841
$ws.Range("O2").Value = 0.2299953477621856
$ws.Range("P2").Value = 0.2299953477621856
$ws.Range("Q2").Value = 10.83258217479144
$ws.Range("R2").Value = 97.49323957312299
$ws.Range("S2").Value = 0.0001581230033989751
$ws.Range("T2").Value = 0.0001581230033989751

# Row 3
$ws.Range("G3").Value = 1.046867666666667
$ws.Range("H3").Value = 3.140603
$ws.Range("I3").Value = 0.000687505225377314
$ws.Range("J3").Value = 0.000687505225377314
$ws.Range("O3").Value = 0.6794731949692173
$ws.Range("P3").Value = 0.6794731949692174
$ws.Range("Q3").Value = 32.00260045121789
$ws.Range("R3").Value = 288.023404060961
$ws.Range("S3").Value = 0.0004671413720451553
$ws.Range("T3").Value = 0.0004671413720451554

# Row 4
$ws.Range("G4").Value = 1.046867666666667
$ws.Range("H4").Value = 3.140603
$ws.Range("I4").Value = 0.000687505225377314
$ws.Range("J4").Value = 0.000687505225377314
$ws.Range("M4").Value = 4.073058666666666
$ws.Range("N4").Value = 12.219176
$ws.Range("O4").Value = 0.09053145726859702
$ws.Range("P4").Value = 0.09053145726859703
$ws.Range("Q4").Value = 4.263953422569778
$ws.Range("R4").Value = 38.375580803128
$ws.Range("S4").Value = 0.00006224084993318346
$ws.Range("T4").Value = 0.00006224084993318347

# Row 5
$ws.Range("H5").Value = 4442.55542
$ws.Range("I5").Value = 0.9725138978974124
$ws.Range("J5").Value = 0.9725138978974125
$ws.Range("M5").Value = 10.34761366666667
$ws.Range("N5").Value = 31.042841
$ws.Range("O5").Value = 0.2299953477621856
$ws.Range("P5").Value = 0.2299953477621856
$ws.Range("Q5").Value = 15323.28239297202
$ws.Range("R5").Value = 137909.5415367482
$ws.Range("S5").Value = 0.223673672150474
$ws.Range("T5").Value = 0.2236736721504741

# Row 6
$ws.Range("H6").Value = 4442.55542
$ws.Range("I6").Value = 0.9725138978974124
$ws.Range("J6").Value = 0.9725138978974125
$ws.Range("O6").Value = 0.6794731949692173
$ws.Range("P6").Value = 0.6794731949692174
$ws.Range("Q6").Value = 45269.43586586795
$ws.Range("R6").Value = 407424.9227928115
$ws.Range("S6").Value = 0.660797125356322
$ws.Range("T6").Value = 0.6607971253563222

# Row 7
$ws.Range("H7").Value = 4442.55542
$ws.Range("I7").Value = 0.9725138978974124
$ws.Range("J7").Value = 0.9725138978974125
$ws.Range("M7").Value = 4.073058666666666
$ws.Range("N7").Value = 12.219176
$ws.Range("O7").Value = 0.09053145726859702
$ws.Range("P7").Value = 0.09053145726859703
$ws.Range("Q7").Value = 6031.596285192657
$ws.Range("R7").Value = 54284.36656673392
$ws.Range("S7").Value = 0.08804310039061632
$ws.Range("T7").Value = 0.08804310039061633

# Row 8
$ws.Range("G8").Value = 40.80635833333333
$ws.Range("H8").Value = 122.419075
$ws.Range("I8").Value = 0.02679859687721029
$ws.Range("J8").Value = 0.0267985968772103
$ws.Range("M8").Value = 10.34761366666667
$ws.Range("N8").Value = 31.042841
$ws.Range("O8").Value = 0.2299953477621856
$ws.Range("P8").Value = 0.2299953477621856
$ws.Range("Q8").Value = 422.2484311768972
$ws.Range("R8").Value = 3800.235880592075
$ws.Range("S8").Value = 0.006163552608312602
$ws.Range("T8").Value = 0.006163552608312604

# Row 9
$ws.Range("G9").Value = 40.80635833333333
$ws.Range("H9").Value = 122.419075
$ws.Range("I9").Value = 0.02679859687721029
$ws.Range("J9").Value = 0.0267985968772103
$ws.Range("O9").Value = 0.6794731949692173
$ws.Range("P9").Value = 0.6794731949692174
$ws.Range("Q9").Value = 1247.444756574669
$ws.Range("R9").Value = 11227.00280917203
$ws.Range("S9").Value = 0.01820892824085017
$ws.Range("T9").Value = 0.01820892824085018

# Row 10
$ws.Range("G10").Value = 40.80635833333333
$ws.Range("H10").Value = 122.419075
$ws.Range("I10").Value = 0.02679859687721029
$ws.Range("J10").Value = 0.0267985968772103
$ws.Range("M10").Value = 4.073058666666666
$ws.Range("N10").Value = 12.219176
$ws.Range("O10").Value = 0.09053145726859702
$ws.Range("P10").Value = 0.09053145726859703
$ws.Range("Q10").Value = 166.2066914646888
$ws.Range("R10").Value = 1495.8602231822
$ws.Range("S10").Value = 0.002426116028047521
$ws.Range("T10").Value = 0.002426116028047522
